# Actualización automática 2025-10-16 14:30:09
#
# Updates the monthly sales / compliance figures for client
# "JARAMILLO CARVAJAL NICOLAS ESTEBAN" (advisor HIDALGO HIDALGO PEDRO
# GUSTAVO) in the "240X80 PORCELANATO" and "PORCELANATO" groups, and
# propagates the resulting totals across the three report sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------
$wsGrupo = $wb.Worksheets("VENTAS POR GRUPO")

# 240X80 PORCELANATO sales for this client go from 95.04 to 552.96
$wsGrupo.Range("D11").Value = 552.96

# PORCELANATO sales for this client go from 0 to 4524.27
$wsGrupo.Range("M11").Value = 4524.27

# Compliance counter for the PORCELANATO column now has one advisor
# meeting it (was 0 of 21, now 1 of 21)
$wsGrupo.Range("M23").Value = "1 de 21"

# ---------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------
$wsMensual = $wb.Worksheets("VENTA MENSUAL")

# October ("octubre") sales for this client rise from 95.04 to 5077.23
$wsMensual.Range("F11").Value = 5077.23

# Column total for October rises accordingly
$wsMensual.Range("F23").Value = 5077.23

# ---------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------
$wsCumpl = $wb.Worksheets("CUMPLIMIENTO MENSUAL")

# Slightly widen the VENTA / CUMPLIMIENTO columns to fit the new values
# (ColumnWidth is expressed in "characters"; the saved OOXML <col width>
# is ColumnWidth + 5/6, so back the offset out to land on exactly 13/24)
$wsCumpl.Columns.Item(4).ColumnWidth = 13 - (5/6)
$wsCumpl.Columns.Item(6).ColumnWidth = 24 - (5/6)

# Row 3: 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 552.96
$wsCumpl.Range("E3").Value = 4951.65890386263
$wsCumpl.Range("F3").Value = 0.1004538206290692

# Row 12: PORCELANATO
$wsCumpl.Range("D12").Value = 4524.27
$wsCumpl.Range("E12").Value = 33215.47
$wsCumpl.Range("F12").Value = 0.1198807940913213

# Row 14: TOTAL
$wsCumpl.Range("D14").Value = 5077.23
$wsCumpl.Range("E14").Value = 50347.51147880389
$wsCumpl.Range("F14").Value = 0.09160584000092609
